$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 33.22094666666666
$ws.Cells.Item(2, 8).Value = 99.66283999999999
$ws.Cells.Item(2, 9).Value = 0.04806943331163595
$ws.Cells.Item(2, 10).Value = 0.04806943331163595
$ws.Cells.Item(2, 13).Value = 10.97489366666667
$ws.Cells.Item(2, 14).Value = 32.924681
$ws.Cells.Item(2, 15).Value = 0.3804091425755969
$ws.Cells.Item(2, 16).Value = 0.3804091425755969
$ws.Cells.Item(2, 17).Value = 364.596357172671
$ws.Cells.Item(2, 18).Value = 3281.367214554039
$ws.Cells.Item(2, 19).Value = 0.01828605191017427
$ws.Cells.Item(2, 20).Value = 0.01828605191017427
$ws.Cells.Item(3, 7).Value = 33.22094666666666
$ws.Cells.Item(3, 8).Value = 99.66283999999999
$ws.Cells.Item(3, 9).Value = 0.04806943331163595
$ws.Cells.Item(3, 10).Value = 0.04806943331163595
$ws.Cells.Item(3, 15).Value = 0.34806939594221
$ws.Cells.Item(3, 16).Value = 0.34806939594221
$ws.Cells.Item(3, 17).Value = 333.6009038705022
$ws.Cells.Item(3, 18).Value = 3002.40813483452
$ws.Cells.Item(3, 19).Value = 0.01673149861606548
$ws.Cells.Item(3, 20).Value = 0.01673149861606548
$ws.Cells.Item(4, 7).Value = 33.22094666666666
$ws.Cells.Item(4, 8).Value = 99.66283999999999
$ws.Cells.Item(4, 9).Value = 0.04806943331163595
$ws.Cells.Item(4, 10).Value = 0.04806943331163595
$ws.Cells.Item(4, 13).Value = 5.513093333333334
$ws.Cells.Item(4, 14).Value = 16.53928
$ws.Cells.Item(4, 15).Value = 0.1910935241443256
$ws.Cells.Item(4, 16).Value = 0.1910935241443256
$ws.Cells.Item(4, 17).Value = 183.1501795950222
$ws.Cells.Item(4, 18).Value = 1648.3516163552
$ws.Cells.Item(4, 19).Value = 0.009185757415141155
$ws.Cells.Item(4, 20).Value = 0.009185757415141155
$ws.Cells.Item(5, 7).Value = 33.22094666666666
$ws.Cells.Item(5, 8).Value = 99.66283999999999
$ws.Cells.Item(5, 9).Value = 0.04806943331163595
$ws.Cells.Item(5, 10).Value = 0.04806943331163595
$ws.Cells.Item(5, 13).Value = 2.320365
$ws.Cells.Item(5, 14).Value = 6.961094999999999
$ws.Cells.Item(5, 15).Value = 0.08042793733786743
$ws.Cells.Item(5, 16).Value = 0.08042793733786745
$ws.Cells.Item(5, 17).Value = 77.08472191219998
$ws.Cells.Item(5, 18).Value = 693.7624972097999
$ws.Cells.Item(5, 19).Value = 0.003866125370255054
$ws.Cells.Item(5, 20).Value = 0.003866125370255055
$ws.Cells.Item(6, 7).Value = 614.8671876666667
$ws.Cells.Item(6, 9).Value = 0.8896891942791112
$ws.Cells.Item(6, 10).Value = 0.8896891942791111
$ws.Cells.Item(6, 13).Value = 10.97489366666667
$ws.Cells.Item(6, 14).Value = 32.924681
$ws.Cells.Item(6, 15).Value = 0.3804091425755969
$ws.Cells.Item(6, 16).Value = 0.3804091425755969
$ws.Cells.Item(6, 17).Value = 6748.102003764046
$ws.Cells.Item(6, 18).Value = 60732.91803387641
$ws.Cells.Item(6, 19).Value = 0.3384459035544903
$ws.Cells.Item(6, 20).Value = 0.3384459035544903
$ws.Cells.Item(7, 7).Value = 614.8671876666667
$ws.Cells.Item(7, 9).Value = 0.8896891942791112
$ws.Cells.Item(7, 10).Value = 0.8896891942791111
$ws.Cells.Item(7, 15).Value = 0.34806939594221
$ws.Cells.Item(7, 16).Value = 0.34806939594221
$ws.Cells.Item(7, 18).Value = 55569.82661019565
$ws.Cells.Item(7, 19).Value = 0.3096735804290418
$ws.Cells.Item(7, 20).Value = 0.3096735804290418
$ws.Cells.Item(8, 7).Value = 614.8671876666667
$ws.Cells.Item(8, 9).Value = 0.8896891942791112
$ws.Cells.Item(8, 10).Value = 0.8896891942791111
$ws.Cells.Item(8, 13).Value = 5.513093333333334
$ws.Cells.Item(8, 14).Value = 16.53928
$ws.Cells.Item(8, 15).Value = 0.1910935241443256
$ws.Cells.Item(8, 16).Value = 0.1910935241443256
$ws.Cells.Item(8, 17).Value = 3389.820193210516
$ws.Cells.Item(8, 18).Value = 30508.38173889464
$ws.Cells.Item(8, 19).Value = 0.170013843527921
$ws.Cells.Item(8, 20).Value = 0.1700138435279209
$ws.Cells.Item(9, 7).Value = 614.8671876666667
$ws.Cells.Item(9, 9).Value = 0.8896891942791112
$ws.Cells.Item(9, 10).Value = 0.8896891942791111
$ws.Cells.Item(9, 13).Value = 2.320365
$ws.Cells.Item(9, 14).Value = 6.961094999999999
$ws.Cells.Item(9, 15).Value = 0.08042793733786743
$ws.Cells.Item(9, 16).Value = 0.08042793733786745
$ws.Cells.Item(9, 17).Value = 1426.716301910165
$ws.Cells.Item(9, 18).Value = 12840.44671719149
$ws.Cells.Item(9, 19).Value = 0.07155586676765811
$ws.Cells.Item(9, 20).Value = 0.07155586676765813
$ws.Cells.Item(10, 7).Value = 42.195614
$ws.Cells.Item(10, 8).Value = 126.586842
$ws.Cells.Item(10, 9).Value = 0.06105543209133513
$ws.Cells.Item(10, 10).Value = 0.06105543209133512
$ws.Cells.Item(10, 13).Value = 10.97489366666667
$ws.Cells.Item(10, 14).Value = 32.924681
$ws.Cells.Item(10, 15).Value = 0.3804091425755969
$ws.Cells.Item(10, 16).Value = 0.3804091425755969
$ws.Cells.Item(10, 17).Value = 463.0923768497113
$ws.Cells.Item(10, 18).Value = 4167.831391647402
$ws.Cells.Item(10, 19).Value = 0.02322604457144738
$ws.Cells.Item(10, 20).Value = 0.02322604457144738
$ws.Cells.Item(11, 7).Value = 42.195614
$ws.Cells.Item(11, 8).Value = 126.586842
$ws.Cells.Item(11, 9).Value = 0.06105543209133513
$ws.Cells.Item(11, 10).Value = 0.06105543209133512
$ws.Cells.Item(11, 15).Value = 0.34806939594221
$ws.Cells.Item(11, 16).Value = 0.34806939594221
$ws.Cells.Item(11, 17).Value = 423.7234751619807
$ws.Cells.Item(11, 18).Value = 3813.511276457826
$ws.Cells.Item(11, 19).Value = 0.02125152736702164
$ws.Cells.Item(11, 20).Value = 0.02125152736702164
$ws.Cells.Item(12, 7).Value = 42.195614
$ws.Cells.Item(12, 8).Value = 126.586842
$ws.Cells.Item(12, 9).Value = 0.06105543209133513
$ws.Cells.Item(12, 10).Value = 0.06105543209133512
$ws.Cells.Item(12, 13).Value = 5.513093333333334
$ws.Cells.Item(12, 14).Value = 16.53928
$ws.Cells.Item(12, 15).Value = 0.1910935241443256
$ws.Cells.Item(12, 16).Value = 0.1910935241443256
$ws.Cells.Item(12, 17).Value = 232.6283582393067
$ws.Cells.Item(12, 18).Value = 2093.65522415376
$ws.Cells.Item(12, 19).Value = 0.01166729768648778
$ws.Cells.Item(12, 20).Value = 0.01166729768648778
$ws.Cells.Item(13, 7).Value = 42.195614
$ws.Cells.Item(13, 8).Value = 126.586842
$ws.Cells.Item(13, 9).Value = 0.06105543209133513
$ws.Cells.Item(13, 10).Value = 0.06105543209133512
$ws.Cells.Item(13, 13).Value = 2.320365
$ws.Cells.Item(13, 14).Value = 6.961094999999999
$ws.Cells.Item(13, 15).Value = 0.08042793733786743
$ws.Cells.Item(13, 16).Value = 0.08042793733786745
$ws.Cells.Item(13, 17).Value = 97.90922587910998
$ws.Cells.Item(13, 18).Value = 881.1830329119898
$ws.Cells.Item(13, 19).Value = 0.004910562466378322
$ws.Cells.Item(13, 20).Value = 0.004910562466378323
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 0.8196073333333334
$ws.Cells.Item(14, 8).Value = 2.458822
$ws.Cells.Item(14, 9).Value = 0.001185940317917725
$ws.Cells.Item(14, 10).Value = 0.001185940317917725
$ws.Cells.Item(14, 13).Value = 10.97489366666667
$ws.Cells.Item(14, 14).Value = 32.924681
$ws.Cells.Item(14, 15).Value = 0.3804091425755969
$ws.Cells.Item(14, 16).Value = 0.3804091425755969
$ws.Cells.Item(14, 17).Value = 8.995103331753556
$ws.Cells.Item(14, 18).Value = 80.95592998578201
$ws.Cells.Item(14, 19).Value = 0.0004511425394849126
$ws.Cells.Item(14, 20).Value = 0.0004511425394849126
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 0.8196073333333334
$ws.Cells.Item(15, 8).Value = 2.458822
$ws.Cells.Item(15, 9).Value = 0.001185940317917725
$ws.Cells.Item(15, 10).Value = 0.001185940317917725
$ws.Cells.Item(15, 15).Value = 0.34806939594221
$ws.Cells.Item(15, 16).Value = 0.34806939594221
$ws.Cells.Item(15, 17).Value = 8.230402040085112
$ws.Cells.Item(15, 18).Value = 74.07361836076601
$ws.Cells.Item(15, 19).Value = 0.0004127895300811351
$ws.Cells.Item(15, 20).Value = 0.0004127895300811351
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 0.8196073333333334
$ws.Cells.Item(16, 8).Value = 2.458822
$ws.Cells.Item(16, 9).Value = 0.001185940317917725
$ws.Cells.Item(16, 10).Value = 0.001185940317917725
$ws.Cells.Item(16, 13).Value = 5.513093333333334
$ws.Cells.Item(16, 14).Value = 16.53928
$ws.Cells.Item(16, 15).Value = 0.1910935241443256
$ws.Cells.Item(16, 16).Value = 0.1910935241443256
$ws.Cells.Item(16, 17).Value = 4.518571725351111
$ws.Cells.Item(16, 18).Value = 40.66714552816001
$ws.Cells.Item(16, 19).Value = 0.00022662551477574
$ws.Cells.Item(16, 20).Value = 0.00022662551477574
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 0.8196073333333334
$ws.Cells.Item(17, 8).Value = 2.458822
$ws.Cells.Item(17, 9).Value = 0.001185940317917725
$ws.Cells.Item(17, 10).Value = 0.001185940317917725
$ws.Cells.Item(17, 13).Value = 2.320365
$ws.Cells.Item(17, 14).Value = 6.961094999999999
$ws.Cells.Item(17, 15).Value = 0.08042793733786743
$ws.Cells.Item(17, 16).Value = 0.08042793733786745
$ws.Cells.Item(17, 17).Value = 1.90178817001
$ws.Cells.Item(17, 18).Value = 17.11609353009
$ws.Cells.Item(17, 19).Value = 0.00009538273357593737
$ws.Cells.Item(17, 20).Value = 0.00009538273357593738
